$wb = $excel.ActiveWorkbook

# ---------- Sheet: Overview ----------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "2593015e-e604-4d8a-a976-6453ef59a653.md"
$ws1.Range("B2").Value = "Handback transform failed"
$ws1.Range("C2").Value = "Handback transform failed"
$ws1.Range("D2").Value = "2016-48-13 14:48:14"
$ws1.Range("A3").Value = "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-47-13 14:47:32"
$ws1.Range("A4").Value = "24e34a03-7cf2-44d5-9af8-b62940d5446f.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-47-13 14:47:32"

# Rebuild hyperlinks on Overview (values changed, ranges shifted)
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md", "", "", "2593015e-e604-4d8a-a976-6453ef59a653.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/24e34a03-7cf2-44d5-9af8-b62940d5446f.md", "", "", "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/725347b5e115091ddf5d6a0898f4be719561ec71/e2e/2593015e-e604-4d8a-a976-6453ef59a653.md", "", "", "24e34a03-7cf2-44d5-9af8-b62940d5446f.md") | Out-Null

# ---------- Sheet: zh-cn ----------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "2593015e-e604-4d8a-a976-6453ef59a653.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handback transform failed"
$ws2.Range("D2").Value = "2593015e-e604-4d8a-a976-6453ef59a653.75fc1b361e66756fecd256450813964da5ae816f.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-13 14:48:11"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("I2").Value = "Include"
$ws2.Range("K2").Value = "Handback file name: wlhuvzrs.tve is different with handoff file name: 2593015e-e604-4d8a-a976-6453ef59a653.75fc1b361e66756fecd256450813964da5ae816f.zh-cn."
$ws2.Range("A3").Value = "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.f0d97fa77fec1060f1a4f1255a3537c7d0c34771.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-13 14:45:37"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"
$ws2.Range("A4").Value = "24e34a03-7cf2-44d5-9af8-b62940d5446f.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "24e34a03-7cf2-44d5-9af8-b62940d5446f.0bf16684cc508d7aa53d6c3e3015f8d97038a05f.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-13 14:45:37"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("I4").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md", "", "", "2593015e-e604-4d8a-a976-6453ef59a653.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b1cb6875ac11b04c2ae16426321b369a9d2d2a3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/101bbb3b-60d8-4015-a9fb-44c5b263c6b0.f0d97fa77fec1060f1a4f1255a3537c7d0c34771.zh-cn.xlf", "", "", "2593015e-e604-4d8a-a976-6453ef59a653.75fc1b361e66756fecd256450813964da5ae816f.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/24e34a03-7cf2-44d5-9af8-b62940d5446f.md", "", "", "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/24e34a03-7cf2-44d5-9af8-b62940d5446f.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b1cb6875ac11b04c2ae16426321b369a9d2d2a3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/24e34a03-7cf2-44d5-9af8-b62940d5446f.0bf16684cc508d7aa53d6c3e3015f8d97038a05f.zh-cn.xlf", "", "", "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.f0d97fa77fec1060f1a4f1255a3537c7d0c34771.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/725347b5e115091ddf5d6a0898f4be719561ec71/e2e/2593015e-e604-4d8a-a976-6453ef59a653.md", "", "", "24e34a03-7cf2-44d5-9af8-b62940d5446f.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/725347b5e115091ddf5d6a0898f4be719561ec71/e2e/2593015e-e604-4d8a-a976-6453ef59a653.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8965926ddfec7e2c91aef7548ff856f3aaaa13ef/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2593015e-e604-4d8a-a976-6453ef59a653.75fc1b361e66756fecd256450813964da5ae816f.zh-cn.xlf", "", "", "24e34a03-7cf2-44d5-9af8-b62940d5446f.0bf16684cc508d7aa53d6c3e3015f8d97038a05f.zh-cn.xlf") | Out-Null

# ---------- Sheet: de-de ----------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "2593015e-e604-4d8a-a976-6453ef59a653.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handback transform failed"
$ws3.Range("D2").Value = "2593015e-e604-4d8a-a976-6453ef59a653.75fc1b361e66756fecd256450813964da5ae816f.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-13 14:48:14"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("I2").Value = "Include"
$ws3.Range("K2").Value = "Handback file name: wlhuvzrs.tve is different with handoff file name: 2593015e-e604-4d8a-a976-6453ef59a653.75fc1b361e66756fecd256450813964da5ae816f.de-de."
$ws3.Range("A3").Value = "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.f0d97fa77fec1060f1a4f1255a3537c7d0c34771.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-13 14:47:32"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"
$ws3.Range("A4").Value = "24e34a03-7cf2-44d5-9af8-b62940d5446f.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "24e34a03-7cf2-44d5-9af8-b62940d5446f.0bf16684cc508d7aa53d6c3e3015f8d97038a05f.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-13 14:47:32"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("I4").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md", "", "", "2593015e-e604-4d8a-a976-6453ef59a653.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbc483f6f836dd93f3633908624a49a71bfb0286/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/101bbb3b-60d8-4015-a9fb-44c5b263c6b0.f0d97fa77fec1060f1a4f1255a3537c7d0c34771.de-de.xlf", "", "", "2593015e-e604-4d8a-a976-6453ef59a653.75fc1b361e66756fecd256450813964da5ae816f.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/24e34a03-7cf2-44d5-9af8-b62940d5446f.md", "", "", "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/fad0d8f8388752b4496f5e3e6c36e117587ba2ae/e2e/24e34a03-7cf2-44d5-9af8-b62940d5446f.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbc483f6f836dd93f3633908624a49a71bfb0286/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/24e34a03-7cf2-44d5-9af8-b62940d5446f.0bf16684cc508d7aa53d6c3e3015f8d97038a05f.de-de.xlf", "", "", "101bbb3b-60d8-4015-a9fb-44c5b263c6b0.f0d97fa77fec1060f1a4f1255a3537c7d0c34771.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/725347b5e115091ddf5d6a0898f4be719561ec71/e2e/2593015e-e604-4d8a-a976-6453ef59a653.md", "", "", "24e34a03-7cf2-44d5-9af8-b62940d5446f.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/725347b5e115091ddf5d6a0898f4be719561ec71/e2e/2593015e-e604-4d8a-a976-6453ef59a653.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f0f0ff4980c887826cae5c5f80d913d2af9759c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2593015e-e604-4d8a-a976-6453ef59a653.75fc1b361e66756fecd256450813964da5ae816f.de-de.xlf", "", "", "24e34a03-7cf2-44d5-9af8-b62940d5446f.0bf16684cc508d7aa53d6c3e3015f8d97038a05f.de-de.xlf") | Out-Null
